$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing entries (per .strip() cleanup of names)
$ws.Range("A3").Value = "Genivaldo"
$ws.Range("A4").Value = "José"
$ws.Range("A5").Value = "João"
$ws.Range("B5").Value = 5

# Insert a new row for Marcelo before the Matias row, shifting Matias down
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = "Marcelo"
$ws.Range("B6").Value = 6
